$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="304.23"},
    @{Cell="E2"; Value="2.03%"},
    @{Cell="D3"; Value="31.94"},
    @{Cell="E3"; Value="1.35%"},
    @{Cell="D4"; Value="5.188"},
    @{Cell="E4"; Value="0.35%"},
    @{Cell="D5"; Value="0.07473"},
    @{Cell="E5"; Value="-0.28%"},
    @{Cell="D6"; Value="2.486"},
    @{Cell="E6"; Value="49.14%"},
    @{Cell="D7"; Value="8.011"},
    @{Cell="E7"; Value="2.94%"},
    @{Cell="D8"; Value="0.9182"},
    @{Cell="E8"; Value="-0.82%"},
    @{Cell="D9"; Value="0.1739"},
    @{Cell="E9"; Value="1.91%"},
    @{Cell="D10"; Value="0.07658"},
    @{Cell="E10"; Value="1.79%"},
    @{Cell="D11"; Value="0.08197"},
    @{Cell="E11"; Value="3.13%"},
    @{Cell="E12"; Value="0.70%"},
    @{Cell="D13"; Value="0.09931"},
    @{Cell="E13"; Value="0.50%"},
    @{Cell="D14"; Value="0.001519"},
    @{Cell="E14"; Value="2.02%"},
    @{Cell="D15"; Value="0.006107"},
    @{Cell="E15"; Value="-7.39%"},
    @{Cell="D16"; Value="3.508"},
    @{Cell="E16"; Value="1.93%"},
    @{Cell="D17"; Value="3.865"},
    @{Cell="E17"; Value="1.79%"},
    @{Cell="D18"; Value="2.228"},
    @{Cell="E18"; Value="-0.06%"},
    @{Cell="D19"; Value="0.3260"},
    @{Cell="E19"; Value="-0.98%"},
    @{Cell="D20"; Value="0.1339"},
    @{Cell="E20"; Value="0.27%"},
    @{Cell="D21"; Value="4.659"},
    @{Cell="E21"; Value="2.11%"},
    @{Cell="D22"; Value="0.04612"},
    @{Cell="E22"; Value="-1.10%"},
    @{Cell="D23"; Value="0.1566"},
    @{Cell="E23"; Value="1.08%"},
    @{Cell="D24"; Value="0.001263"},
    @{Cell="E24"; Value="3.55%"},
    @{Cell="D25"; Value="0.004521"},
    @{Cell="E25"; Value="2.39%"},
    @{Cell="D26"; Value="0.0001300"},
    @{Cell="E26"; Value="-7.09%"},
    @{Cell="D27"; Value="0.0002743"},
    @{Cell="E27"; Value="51.65%"},
    @{Cell="D39"; Value="0.01771"},
    @{Cell="E39"; Value="6.69%"},
    @{Cell="D40"; Value="0.04544"},
    @{Cell="E40"; Value="0.20%"},
    @{Cell="D41"; Value="0.007365"},
    @{Cell="E41"; Value="5.27%"},
    @{Cell="D42"; Value="0.1361"},
    @{Cell="E42"; Value="1.41%"},
    @{Cell="D43"; Value="0.002141"},
    @{Cell="E43"; Value="3.95%"},
    @{Cell="E44"; Value="-17.82%"},
    @{Cell="D45"; Value="0.00006483"},
    @{Cell="E45"; Value="6.99%"},
    @{Cell="E46"; Value="-57.48%"}
)

foreach ($chg in $changes) {
    $ws.Range($chg.Cell).Value = "'" + $chg.Value
    $ws.Range($chg.Cell).Style = "Normal"
}
